$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "a"
$ws.Range("B2").Value = "E:\storage\905 - Copy (14).mp4"
$ws.Range("C2").Value = "a"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "11/14/2025"
$ws.Range("E2").ClearFormats()
$ws.Range("G2").Value = "E:\ghep xong\905 - Copy (14).mp4"

# Row 3
$ws.Range("A3").Value = "b"
$ws.Range("B3").Value = "s"
$ws.Range("C3").Value = "a"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "11/14/2025"
$ws.Range("E3").ClearFormats()

# Row 4
$ws.Range("A4").Value = "c"
$ws.Range("B4").Value = "E:\storage\902 - Copy (13).mp4"
$ws.Range("C4").Value = "a"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "11/14/2025"
$ws.Range("E4").ClearFormats()
$ws.Range("G4").Value = "E:\ghep xong\902 - Copy (13).mp4"

# Row 5 (new row)
$ws.Range("A5").Value = "a"
$ws.Range("B5").Value = "s"
$ws.Range("C5").Value = "s"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "11/14/2025"
$ws.Range("E5").ClearFormats()
$ws.Range("G5").Font.Bold = $false
